$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first data point (0.0 min / 20.0 °C) is removed: the profile's first
# point is always room temperature, so it no longer needs to be listed
# explicitly. Deleting row 2 shifts the remaining data rows up by one.
$ws.Rows.Item(2).Delete()

# Add a new "Beschreibung" column describing each remaining phase.
# Copy formatting from the existing header/data cells so the new column
# matches the existing style (s="1") instead of creating a new style.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A2").Copy()
$ws.Range("C2:C7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C1").Value = "Beschreibung"
$ws.Range("C2").Value = "Aufheizen 1"
$ws.Range("C3").Value = "Aufheizen 2"
$ws.Range("C4").Value = "Aufheizen 3"
$ws.Range("C5").Value = "Halten"
$ws.Range("C6").Value = "Abkühlen 1"
$ws.Range("C7").Value = "Abkühlen 2"
